$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from G1 to H1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add new "Save" column values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
